$d = $word.ActiveDocument

# The document currently ends with an empty, list-styled paragraph
# (Paragraphedeliste / numId 2) that follows the "... exercise ?" bullet.
# We need to insert two new bullet paragraphs *before* that trailing
# empty paragraph, each reproducing Word's grammar-check proofErr
# bookends ("gramStart"/"gramEnd") around the final "word ?" run, just
# like the existing bullets in the document.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# --- New paragraph 1 --------------------------------------------------
$lastPara.Range.InsertParagraphBefore()
$target = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range

$xml1 = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragraphedeliste"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:rPr>
                <w:rFonts w:ascii="Century" w:hAnsi="Century"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Century" w:hAnsi="Century"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
              <w:t xml:space="preserve">Is it possible to model a flexible extrusion / flexible section that can adapt to 3D model in an </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Century" w:hAnsi="Century"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
              <w:t>assembly ?</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Century" w:hAnsi="Century"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
              <w:t xml:space="preserve"> (just for curiosity)</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$target.InsertXML($xml1)

# --- New paragraph 2 --------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()
$target2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range

$xml2 = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragraphedeliste"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:rPr>
                <w:rFonts w:ascii="Century" w:hAnsi="Century"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Century" w:hAnsi="Century"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
              <w:t xml:space="preserve">Minor clashes with assembly </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Century" w:hAnsi="Century"/>
                <w:lang w:val="en-GB"/>
              </w:rPr>
              <w:t>strip ?</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$target2.InsertXML($xml2)
